$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.077.74"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.789.32"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4293"
$ws.Range("E7").Value = "  -3.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3620"
$ws.Range("E8").Value = "  -3.21%  "
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07508"
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.114"
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.63"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.140"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.308"
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "1.803.45"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.16"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06350"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.978"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").Value = "28.087.02"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  -2.42%  "
$ws.Range("E25").Value = "  -7.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.84"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").Value = "2.006.63"
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.179"
$ws.Range("E29").Value = "  -7.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.16"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.161"
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.755"
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08991"
$ws.Range("E33").Value = "  -3.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.520"
$ws.Range("E34").Value = "  -3.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.62"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02327"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.087"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6474"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2113"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06059"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.187"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.419"
$ws.Range("E42").Value = "  -0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.880"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.60"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.50"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.983"
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.155"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06951"
$ws.Range("E51").Value = "  +0.50%  "
